$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

# Row 7
$ws.Range("B7").Value = 0.1838799416186027
$ws.Range("C7").Value = 0.7695350798161764
$ws.Range("D7").Value = 1.48474280445927
$ws.Range("E7").Value = 1.218500227517119
$ws.Range("F7").Value = 1.220715086378263
$ws.Range("G7").Value = 38

# Row 8
$ws.Range("B8").Value = 0.2202601043985484
$ws.Range("C8").Value = 0.7025675926365696
$ws.Range("D8").Value = 1.450233190988454
$ws.Range("E8").Value = 1.204256281274237
$ws.Range("F8").Value = 1.200273013662533
$ws.Range("G8").Value = 37

# Row 9
$ws.Range("B9").Value = 0.05458795430075868
$ws.Range("C9").Value = 0.6814427688448653
$ws.Range("D9").Value = 0.8599296228729056
$ws.Range("E9").Value = 0.9273239039693227
$ws.Range("F9").Value = 0.9497643907933797
$ws.Range("G9").Value = 20

# Row 10
$ws.Range("B10").Value = 0.3270187541837685
$ws.Range("C10").Value = 0.5806918712726016
$ws.Range("D10").Value = 0.7433548639178584
$ws.Range("E10").Value = 0.8621802966420994
$ws.Range("F10").Value = 0.8303300939129272
$ws.Range("G10").Value = 13

# Row 11
$ws.Range("B11").Value = 0.2613329855841933
$ws.Range("C11").Value = 0.4854844728135347
$ws.Range("D11").Value = 0.3301213188269617
$ws.Range("E11").Value = 0.5745618494356911
$ws.Range("F11").Value = 0.5720865204151965
$ws.Range("G11").Value = 5
